$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 123
$ws.Range("E2").Value = 123
$ws.Range("F2").Value = 123
$ws.Range("G2").Value = 123
$ws.Range("B3").Value = 246
